$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing A12 timestamp to its slightly adjusted value
$ws.Cells.Item(12, 1).Value = 45877.45849543982

# Add new row 13 with the new sensor reading
$ws.Cells.Item(13, 1).Value = 45877.5001959631
$ws.Cells.Item(13, 2).Value = 2025
$ws.Cells.Item(13, 3).Value = 32
$ws.Cells.Item(13, 4).Value = 16.61
$ws.Cells.Item(13, 5).Value = 85.93000000000001
$ws.Cells.Item(13, 6).Value = 488.67
$ws.Cells.Item(13, 7).Value = 13.47
$ws.Cells.Item(13, 8).Value = "ESE"
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = "12:00:16"

# Match the existing date/time number format used by column A (style index 2 -> numFmt 165)
$ws.Cells.Item(13, 1).NumberFormat = $ws.Cells.Item(12, 1).NumberFormat
